$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 669005.4399999999   # H40: 557771.1 -> 669005.4399999999
$ws.Cells.Item(40, 9).Value = 1666.6666   # I40: 1599.875 -> 1666.6666
$ws.Cells.Item(40, 10).Value = 1113898   # J40: 1002708.1 -> 1113898
$ws.Cells.Item(40, 11).Value = 1666.6666   # K40: 1599.875 -> 1666.6666
$ws.Cells.Item(40, 12).Value = 1113898   # L40: 1002708.1 -> 1113898
$ws.Cells.Item(40, 13).Value = -1491.6666   # M40: -1424.875 -> -1491.6666
$ws.Cells.Item(40, 14).Value = -1114248   # N40: -1003058.1 -> -1114248

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 1043.6072   # H112: 15551.363 -> 1043.6072
$ws.Cells.Item(112, 9).Value = 500   # I112: 475 -> 500
$ws.Cells.Item(112, 10).Value = 1053.491   # J112: 16733.824 -> 1053.491
$ws.Cells.Item(112, 11).Value = 1500   # K112: 1425 -> 1500
$ws.Cells.Item(112, 12).Value = 3160.473   # L112: 50201.472 -> 3160.473
$ws.Cells.Item(112, 13).Value = -392   # M112: -317 -> -392
$ws.Cells.Item(112, 14).Value = -5376.473   # N112: -52417.472 -> -5376.473

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1016.8033   # H137: 1124.1569 -> 1016.8033
$ws.Cells.Item(137, 9).Value = 913.4583   # I137: 1017.1842 -> 913.4583
$ws.Cells.Item(137, 10).Value = 1398.3846   # J137: 1436.8462 -> 1398.3846
$ws.Cells.Item(137, 11).Value = 2740.3749   # K137: 3051.5526 -> 2740.3749
$ws.Cells.Item(137, 12).Value = 4195.1538   # L137: 4310.5386 -> 4195.1538
$ws.Cells.Item(137, 13).Value = -190.3748999999998   # M137: -501.5526 -> -190.3748999999998
$ws.Cells.Item(137, 14).Value = -9295.1538   # N137: -9410.5386 -> -9295.1538

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2359.22   # H32: 10642327 -> 2359.22
$ws.Cells.Item(32, 9).Value = 2130.6067   # I32: 3795.1604 -> 2130.6067
$ws.Cells.Item(32, 10).Value = 4208.909   # J32: 76928560 -> 4208.909
$ws.Cells.Item(32, 11).Value = 2130.6067   # K32: 3795.1604 -> 2130.6067
$ws.Cells.Item(32, 12).Value = 4208.909   # L32: 76928560 -> 4208.909
$ws.Cells.Item(32, 13).Value = -1843.6067   # M32: -3508.1604 -> -1843.6067
$ws.Cells.Item(32, 14).Value = -4782.909   # N32: -76929134 -> -4782.909

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1306.3143   # H61: 1316.2954 -> 1306.3143
$ws.Cells.Item(61, 9).Value = 1000.7241   # I61: 894.6177 -> 1000.7241
$ws.Cells.Item(61, 10).Value = 2783.3333   # J61: 2750 -> 2783.3333
$ws.Cells.Item(61, 11).Value = 1000.7241   # K61: 894.6177 -> 1000.7241
$ws.Cells.Item(61, 12).Value = 2783.3333   # L61: 2750 -> 2783.3333
$ws.Cells.Item(61, 13).Value = -788.7241   # M61: -682.6177 -> -788.7241
$ws.Cells.Item(61, 14).Value = -3207.3333   # N61: -3174 -> -3207.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 823.4897999999999   # H74: 918.4286 -> 823.4897999999999
$ws.Cells.Item(74, 9).Value = 846.8298   # I74: 918.4286 -> 846.8298
$ws.Cells.Item(74, 10).Value = 275   # J74: 0 -> 275
$ws.Cells.Item(74, 11).Value = 846.8298   # K74: 918.4286 -> 846.8298
$ws.Cells.Item(74, 12).Value = 275   # L74: 0 -> 275
$ws.Cells.Item(74, 13).Value = 27.17020000000002   # M74: -44.42859999999996 -> 27.17020000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 823.4897999999999   # H77: 918.4286 -> 823.4897999999999
$ws.Cells.Item(77, 9).Value = 846.8298   # I77: 918.4286 -> 846.8298
$ws.Cells.Item(77, 10).Value = 275   # J77: 0 -> 275
$ws.Cells.Item(77, 11).Value = 4234.148999999999   # K77: 4592.143 -> 4234.148999999999
$ws.Cells.Item(77, 12).Value = 1375   # L77: 0 -> 1375
$ws.Cells.Item(77, 13).Value = 133.8510000000006   # M77: -224.143 -> 133.8510000000006

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1051902.2   # H132: 1132801.2 -> 1051902.2
$ws.Cells.Item(132, 9).Value = 1208.7709   # I132: 1307.9767 -> 1208.7709
$ws.Cells.Item(132, 10).Value = 7356063   # J132: 6538824.5 -> 7356063
$ws.Cells.Item(132, 11).Value = 3626.3127   # K132: 3923.9301 -> 3626.3127
$ws.Cells.Item(132, 12).Value = 22068189   # L132: 19616473.5 -> 22068189
$ws.Cells.Item(132, 13).Value = -1096.3127   # M132: -1393.9301 -> -1096.3127
$ws.Cells.Item(132, 14).Value = -22073249   # N132: -19621533.5 -> -22073249

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1306.3143   # H136: 1316.2954 -> 1306.3143
$ws.Cells.Item(136, 9).Value = 1000.7241   # I136: 894.6177 -> 1000.7241
$ws.Cells.Item(136, 10).Value = 2783.3333   # J136: 2750 -> 2783.3333
$ws.Cells.Item(136, 11).Value = 3002.1723   # K136: 2683.8531 -> 3002.1723
$ws.Cells.Item(136, 12).Value = 8349.999899999999   # L136: 8250 -> 8349.999899999999
$ws.Cells.Item(136, 13).Value = -452.1723000000002   # M136: -133.8531000000003 -> -452.1723000000002
$ws.Cells.Item(136, 14).Value = -13449.9999   # N136: -13350 -> -13449.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 1689387.6   # H22: 1501693.5 -> 1689387.6
$ws.Cells.Item(22, 9).Value = 1689387.6   # I22: 1501693.5 -> 1689387.6
$ws.Cells.Item(22, 11).Value = 1689387.6   # K22: 1501693.5 -> 1689387.6
$ws.Cells.Item(22, 13).Value = -1689214.6   # M22: -1501520.5 -> -1689214.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2139629.8   # H134: 2587316 -> 2139629.8
$ws.Cells.Item(134, 9).Value = 898.34283   # I134: 963.0606 -> 898.34283
$ws.Cells.Item(134, 10).Value = 6542900.5   # J134: 11122281 -> 6542900.5
$ws.Cells.Item(134, 11).Value = 2695.02849   # K134: 2889.1818 -> 2695.02849
$ws.Cells.Item(134, 12).Value = 19628701.5   # L134: 33366843 -> 19628701.5
$ws.Cells.Item(134, 13).Value = -160.0284900000001   # M134: -354.1818000000003 -> -160.0284900000001
$ws.Cells.Item(134, 14).Value = -19633771.5   # N134: -33371913 -> -19633771.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1199.3158   # H31: 1272.2354 -> 1199.3158
$ws.Cells.Item(31, 9).Value = 890.9375   # I31: 1026.6522 -> 890.9375
$ws.Cells.Item(31, 10).Value = 1594.04   # J31: 1473.9642 -> 1594.04
$ws.Cells.Item(31, 11).Value = 890.9375   # K31: 1026.6522 -> 890.9375
$ws.Cells.Item(31, 12).Value = 1594.04   # L31: 1473.9642 -> 1594.04
$ws.Cells.Item(31, 13).Value = -595.9375   # M31: -731.6522 -> -595.9375
$ws.Cells.Item(31, 14).Value = -2184.04   # N31: -2063.9642 -> -2184.04

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 1199.3158   # H34: 1272.2354 -> 1199.3158
$ws.Cells.Item(34, 9).Value = 890.9375   # I34: 1026.6522 -> 890.9375
$ws.Cells.Item(34, 10).Value = 1594.04   # J34: 1473.9642 -> 1594.04
$ws.Cells.Item(34, 11).Value = 890.9375   # K34: 1026.6522 -> 890.9375
$ws.Cells.Item(34, 12).Value = 1594.04   # L34: 1473.9642 -> 1594.04
$ws.Cells.Item(34, 13).Value = -688.9375   # M34: -824.6522 -> -688.9375
$ws.Cells.Item(34, 14).Value = -1998.04   # N34: -1877.9642 -> -1998.04

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 14706854   # H58: 12821325 -> 14706854
$ws.Cells.Item(58, 9).Value = 18519326   # I58: 17857926 -> 18519326
$ws.Cells.Item(58, 10).Value = 1603.8572   # J58: 883.9545000000001 -> 1603.8572
$ws.Cells.Item(58, 11).Value = 18519326   # K58: 17857926 -> 18519326
$ws.Cells.Item(58, 12).Value = 1603.8572   # L58: 883.9545000000001 -> 1603.8572
$ws.Cells.Item(58, 13).Value = -18519123   # M58: -17857723 -> -18519123
$ws.Cells.Item(58, 14).Value = -2009.8572   # N58: -1289.9545 -> -2009.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 7247842.5   # H132: 6804343.5 -> 7247842.5
$ws.Cells.Item(132, 9).Value = 1110.6875   # I132: 1111.9375 -> 1110.6875
$ws.Cells.Item(132, 10).Value = 23811802   # J132: 19610426 -> 23811802
$ws.Cells.Item(132, 11).Value = 3332.0625   # K132: 3335.8125 -> 3332.0625
$ws.Cells.Item(132, 12).Value = 71435406   # L132: 58831278 -> 71435406
$ws.Cells.Item(132, 13).Value = -802.0625   # M132: -805.8125 -> -802.0625
$ws.Cells.Item(132, 14).Value = -71440466   # N132: -58836338 -> -71440466

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 13158807   # H134: 11111913 -> 13158807
$ws.Cells.Item(134, 9).Value = 960.2759   # I134: 848.9706 -> 960.2759
$ws.Cells.Item(134, 10).Value = 55556310   # J134: 45455200 -> 55556310
$ws.Cells.Item(134, 11).Value = 2880.8277   # K134: 2546.9118 -> 2880.8277
$ws.Cells.Item(134, 12).Value = 166668930   # L134: 136365600 -> 166668930
$ws.Cells.Item(134, 13).Value = -345.8276999999998   # M134: -11.91179999999986 -> -345.8276999999998
$ws.Cells.Item(134, 14).Value = -166674000   # N134: -136370670 -> -166674000

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 14706854   # H136: 12821325 -> 14706854
$ws.Cells.Item(136, 9).Value = 18519326   # I136: 17857926 -> 18519326
$ws.Cells.Item(136, 10).Value = 1603.8572   # J136: 883.9545000000001 -> 1603.8572
$ws.Cells.Item(136, 11).Value = 55557978   # K136: 53573778 -> 55557978
$ws.Cells.Item(136, 12).Value = 4811.571599999999   # L136: 2651.8635 -> 4811.571599999999
$ws.Cells.Item(136, 13).Value = -55555428   # M136: -53571228 -> -55555428
$ws.Cells.Item(136, 14).Value = -9911.571599999999   # N136: -7751.8635 -> -9911.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 25.105263   # H2: 42.692307 -> 25.105263
$ws.Cells.Item(2, 9).Value = 16   # I2: 18.571428 -> 16
$ws.Cells.Item(2, 10).Value = 35.22222   # J2: 70.833336 -> 35.22222
$ws.Cells.Item(2, 11).Value = 16   # K2: 18.571428 -> 16
$ws.Cells.Item(2, 12).Value = 35.22222   # L2: 70.833336 -> 35.22222
$ws.Cells.Item(2, 13).Value = 97   # M2: 94.428572 -> 97
$ws.Cells.Item(2, 14).Value = -261.22222   # N2: -296.833336 -> -261.22222

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1758.75   # H61: 2025.0625 -> 1758.75
$ws.Cells.Item(61, 9).Value = 1511.8334   # I61: 1984.6923 -> 1511.8334
$ws.Cells.Item(61, 10).Value = 2499.5   # J61: 2200 -> 2499.5
$ws.Cells.Item(61, 11).Value = 1511.8334   # K61: 1984.6923 -> 1511.8334
$ws.Cells.Item(61, 12).Value = 2499.5   # L61: 2200 -> 2499.5
$ws.Cells.Item(61, 13).Value = -1309.8334   # M61: -1782.6923 -> -1309.8334
$ws.Cells.Item(61, 14).Value = -2903.5   # N61: -2604 -> -2903.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 1758.75   # H113: 2025.0625 -> 1758.75
$ws.Cells.Item(113, 9).Value = 1511.8334   # I113: 1984.6923 -> 1511.8334
$ws.Cells.Item(113, 10).Value = 2499.5   # J113: 2200 -> 2499.5
$ws.Cells.Item(113, 11).Value = 1511.8334   # K113: 1984.6923 -> 1511.8334
$ws.Cells.Item(113, 12).Value = 2499.5   # L113: 2200 -> 2499.5
$ws.Cells.Item(113, 13).Value = 658.1666   # M113: 185.3077000000001 -> 658.1666
$ws.Cells.Item(113, 14).Value = -6839.5   # N113: -6540 -> -6839.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 14930058   # H132: 13162085 -> 14930058
$ws.Cells.Item(132, 9).Value = 22728532   # I132: 19231892 -> 22728532
$ws.Cells.Item(132, 10).Value = 11237.305   # J132: 10836.167 -> 11237.305
$ws.Cells.Item(132, 11).Value = 68185596   # K132: 57695676 -> 68185596
$ws.Cells.Item(132, 12).Value = 33711.915   # L132: 32508.501 -> 33711.915
$ws.Cells.Item(132, 13).Value = -68183066   # M132: -57693146 -> -68183066
$ws.Cells.Item(132, 14).Value = -38771.915   # N132: -37568.501 -> -38771.915

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 32092640   # H136: 34331604 -> 32092640
$ws.Cells.Item(136, 9).Value = 11339382   # I136: 11615951 -> 11339382
$ws.Cells.Item(136, 10).Value = 250001860   # J136: 500002500 -> 250001860
$ws.Cells.Item(136, 11).Value = 34018146   # K136: 34847853 -> 34018146
$ws.Cells.Item(136, 12).Value = 750005580   # L136: 1500007500 -> 750005580
$ws.Cells.Item(136, 13).Value = -34015596   # M136: -34845303 -> -34015596
$ws.Cells.Item(136, 14).Value = -750010680   # N136: -1500012600 -> -750010680

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(124, 8).Value = 24214.5   # H124: 28964.5 -> 24214.5
$ws.Cells.Item(124, 10).Value = 24214.5   # J124: 28964.5 -> 24214.5
$ws.Cells.Item(124, 12).Value = 24214.5   # L124: 28964.5 -> 24214.5
$ws.Cells.Item(124, 14).Value = -34034.5   # N124: -38784.5 -> -34034.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 18641.06   # H132: 19812 -> 18641.06
$ws.Cells.Item(132, 9).Value = 22173.918   # I132: 24101.244 -> 22173.918
$ws.Cells.Item(132, 11).Value = 66521.754   # K132: 72303.73199999999 -> 66521.754
$ws.Cells.Item(132, 13).Value = -63991.754   # M132: -69773.73199999999 -> -63991.754

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 9260394   # H136: 10001241 -> 9260394
$ws.Cells.Item(136, 9).Value = 11628536   # I136: 13889618 -> 11628536
$ws.Cells.Item(136, 10).Value = 3110.818   # J136: 2558.5715 -> 3110.818
$ws.Cells.Item(136, 11).Value = 34885608   # K136: 41668854 -> 34885608
$ws.Cells.Item(136, 12).Value = 9332.454000000002   # L136: 7675.7145 -> 9332.454000000002
$ws.Cells.Item(136, 13).Value = -34883058   # M136: -41666304 -> -34883058
$ws.Cells.Item(136, 14).Value = -14432.454   # N136: -12775.7145 -> -14432.454
